$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 289 (before the current row 290),
# shifting the existing rows 290:343 down to 292:345.
$ws.Rows("290:291").Insert()

# Populate the newly inserted rows with this week's price data
# (same market/category/unit/origin/classification as the rest of the sheet).

# Row 290 - Primera
$ws.Range("A290").Value = 8
$ws.Range("B290").Value = "Terminal La Palmera de La Serena"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value = 44617
$ws.Range("E290").Value = 4
$ws.Range("F290").Value = 100112009
$ws.Range("G290").Value = "Acelga"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 2460
$ws.Range("K290").Value = 500
$ws.Range("L290").Value = 600
$ws.Range("M290").Value = 550
$ws.Range("N290").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O290").Value = "Provincia del Elquí"
$ws.Range("P290").Value = 275
$ws.Range("Q290").Value = 2
$ws.Range("R290").Value = "Hortaliza"

# Row 291 - Segunda
$ws.Range("A291").Value = 8
$ws.Range("B291").Value = "Terminal La Palmera de La Serena"
$ws.Range("C291").Value = "Coquimbo"
$ws.Range("D291").Value = 44617
$ws.Range("E291").Value = 4
$ws.Range("F291").Value = 100112009
$ws.Range("G291").Value = "Acelga"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Segunda"
$ws.Range("J291").Value = 1340
$ws.Range("K291").Value = 400
$ws.Range("L291").Value = 450
$ws.Range("M291").Value = 425
$ws.Range("N291").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O291").Value = "Provincia del Elquí"
$ws.Range("P291").Value = 212
$ws.Range("Q291").Value = 2
$ws.Range("R291").Value = "Hortaliza"
